$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order of assignment matters: it determines how new entries are appended to the
# shared strings table (after orphaned entries for old values are pruned).
# Target shared-string order requires: C5, D5, G5, H5(->I5), D4, D3 style sequence
# matching the underlying edits: C5, D5, D4, D3 introduce the brand-new strings;
# other edits reuse already-existing strings.

# --- Row 5 updates (previously blank cells now filled in) ---
$ws.Range("C5").Value = "Select the language to be utilized"
$ws.Range("D5").Value = "Walker (Mentor)"
$ws.Range("E5").Value = (Get-Date -Year 2025 -Month 2 -Day 7).Date
$ws.Range("F5").Value = (Get-Date -Year 2025 -Month 2 -Day 7).Date
$ws.Range("G5").Value = "High"
$ws.Range("H5").Value = "Complete"
$ws.Range("I5").Value = "Yes"

# --- Row 4 updates ---
# D4: "Donovan Ester" -> "Everyone"
$ws.Range("D4").Value = "Everyone"
# C4 stays "Set Up environment" (text unchanged)
$ws.Range("C4").Value = "Set Up environment"
# F4: finish date 2/4/2025 -> 2/16/2025
$ws.Range("F4").Value = (Get-Date -Year 2025 -Month 2 -Day 16).Date
# H4: "Not Started" -> "In Progress"
$ws.Range("H4").Value = "In Progress"

# --- Row 3 updates ---
# D3: "Kameron Smith" -> "Kameron Smith and  Walker (Mentor)"
$ws.Range("D3").Value = "Kameron Smith and  Walker (Mentor)"
# H3: "In Progress" -> "Complete"
$ws.Range("H3").Value = "Complete"

# --- Selection change: active cell moves from N16 to H6 ---
$ws.Range("H6").Select()

$wb.Save()
